$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("user")

# Row 2 - Rahi Shah
$ws.Range("A2").Value = "Rahi"
$ws.Range("B2").Value = "Shah"
$ws.Range("C2").Value = "rahiy@poonamcoatings.com"
$ws.Range("D2").Value = 9838981681
$ws.Range("F2").Value = [DateTime]"1981-12-13"

# Row 3 - Trish Patdel
$ws.Range("A3").Value = "Trish"
$ws.Range("B3").Value = "Patdel"
$ws.Range("C3").Value = "trish@poonamcoatings.com"
$ws.Range("D3").Value = 9827517841
$ws.Range("F3").Value = [DateTime]"1991-11-28"

# Row 4 - Krish Joshi
$ws.Range("A4").Value = "Krish"
$ws.Range("B4").Value = "Joshi"
$ws.Range("C4").Value = "krish@poonamcoatings.com"
$ws.Range("D4").Value = 9875678761
$ws.Range("F4").Value = [DateTime]"1981-10-22"
$ws.Range("G4").Value = "Admin"

# Row 5 - Nrish Kashdyap
$ws.Range("A5").Value = "Nrish"
$ws.Range("B5").Value = "Kashdyap"
$ws.Range("C5").Value = "nrish@poonamcoatings.com"
$ws.Range("D5").Value = 9809167871
$ws.Range("E5").Value = "Other"
$ws.Range("F5").Value = [DateTime]"1987-12-26"

$ws.Range("B9").Select()
